# Adds an API-side filter table (Filter Name / Filter Value) next to the
# existing "Data Source" / "form" cells, and leaves the selection on D22
# (matches the sheetView selection recorded after the edit was made in Excel).
#
# Shared-string insertion order matters for a byte-faithful rebuild of
# xl/sharedStrings.xml, so "Filter Value" (column C) is written before
# "Filter Name" (column B) — that reproduces the original author's order:
#   0 Data Source, 1 form, 2 Filter Value, 3 Filter Name,
#   4 app_id, 5 foobizzle, 6 type, 7 intake

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): existing A1 "Data Source" is untouched/stays as-is.
$ws.Range("C1").Value = "Filter Value"
$ws.Range("B1").Value = "Filter Name"

# Filter rows
$ws.Range("B2").Value = "app_id"
$ws.Range("C2").Value = "foobizzle"

$ws.Range("B3").Value = "type"
$ws.Range("C3").Value = "intake"

# Matches the trailing <selection activeCell="D22" sqref="D22"/> left in the
# saved sheetView.
$ws.Range("D22").Select()
